$wb = $excel.ActiveWorkbook

# Sheets that contain the "展览" (exhibition) rows affected by this update:
# 展览 (rId1 / sheet1) and 全部类型 (rId4 / sheet4)
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5170
    $ws.Range("F3").Value = 158
    $ws.Range("F4").Value = 906
}
